# The author regrouped the weekly price records for
# "Hortaliza, Vega Monumental Concepcion - Sandia": every 3-row (Extra /
# Primera / Segunda) block of columns A:R moves as a whole unit to a new
# row position (the rows are re-sorted into a different weekly order); no
# individual cell value is edited "in place" - this was verified by
# diffing the original workbook against the target OOXML, which shows that
# for every resulting row, the full set of columns A:R matches exactly one
# of the original rows.
#
# $permutation[i-1] gives the 1-based row offset (within A2:R59, so 1 ==
# sheet row 2) of the ORIGINAL row whose contents should end up at new
# offset i (new sheet row i+1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("A2:R59")

# Read the original values twice: $target will be overwritten in place and
# written back to the sheet, $source keeps the untouched original data to
# copy rows from (the COM value arrays returned here are 1-based).
$target = $range.Value()
$source = $range.Value()

$permutation = @(21,22,23,41,42,43,38,39,40,56,57,58,35,36,37,29,30,31,47,48,49,3,4,5,1,2,9,10,11,18,19,20,53,54,55,27,28,12,13,14,6,7,8,50,51,52,24,25,26,44,45,46,15,16,17,32,33,34)

$rowCount = $target.GetLength(0)
$colCount = $target.GetLength(1)

for ($i = 1; $i -le $rowCount; $i++) {
    $srcRow = $permutation[$i - 1]
    for ($j = 1; $j -le $colCount; $j++) {
        $target[$i, $j] = $source[$srcRow, $j]
    }
}

$range.Value = $target
